# The commit adds one new weekly price record for "Zanahoria" at
# "Terminal La Palmera de La Serena", inserted as a new row at sheet
# position 140 (pushing all the existing rows from 140 down to 141+).
#
# The new row reuses the same Volumen/Precio/Unidad/Origen/etc. values
# that (old) row 140 had, only the Fecha (date, column D) is new.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 140; everything at/after row 140
# (old rows 140:273) shifts down to 141:274.
$ws.Rows(140).Insert()

# The row that is now at 141 holds what used to be row 140's data.
# Duplicate that whole row into the freshly inserted row 140.
$src = $ws.Range("A141:R141")
$dst = $ws.Range("A140:R140")
$dst.Value2 = $src.Value2

# Give the new row its own date (Fecha, column D) - 2022-01-25.
$ws.Range("D140").Value2 = 44586
